$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-8 from 45221 (2023-10-22)
# to 45224 (2023-10-25), keeping existing number formatting intact.
$ws.Range("C2:C8").Value = 45224
